# Applies the 2025-12-19 daily update to the violent-crime-full-year workbook.
# For each affected sheet, update the specific cells in column K/L (year 2025/2024 totals)
# with the new cumulative counts.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("L2").Value = 6399
$ws.Range("L3").Value = 6895
$ws.Range("L4").Value = 1708
$ws.Range("L5").Value = 406
$ws.Range("L6").Value = 5668
$ws.Range("L7").Value = 21076

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("L8").Value = 1393
$ws.Range("L11").Value = 349
$ws.Range("L12").Value = 49
$ws.Range("L15").Value = 176
$ws.Range("L19").Value = 574
$ws.Range("L20").Value = 532
$ws.Range("L24").Value = 64
$ws.Range("L25").Value = 127
$ws.Range("L27").Value = 182
$ws.Range("L29").Value = 1180
$ws.Range("L31").Value = 210
$ws.Range("L33").Value = 950
$ws.Range("K37").Value = 909
$ws.Range("L42").Value = 666
$ws.Range("L48").Value = 275
$ws.Range("L53").Value = 237
$ws.Range("L54").Value = 456
$ws.Range("L60").Value = 143
$ws.Range("K63").Value = 180
$ws.Range("L63").Value = 62
$ws.Range("L65").Value = 415
$ws.Range("L67").Value = 729
$ws.Range("L69").Value = 45
$ws.Range("L77").Value = 141
$ws.Range("L79").Value = 584
$ws.Range("L82").Value = 30
$ws.Range("L85").Value = 1046
$ws.Range("L89").Value = 285
$ws.Range("L90").Value = 223
$ws.Range("L94").Value = 257
$ws.Range("L95").Value = 296
$ws.Range("L99").Value = 364
$ws.Range("L101").Value = 21076

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("L3").Value = 99
$ws.Range("L7").Value = 349

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("L3").Value = 84
$ws.Range("L7").Value = 285

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("L2").Value = 315
$ws.Range("L7").Value = 1046

$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range("L2").Value = 19
$ws.Range("L7").Value = 45

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("L4").Value = 30
$ws.Range("L7").Value = 237

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("L2").Value = 426
$ws.Range("L6").Value = 337
$ws.Range("L7").Value = 1393

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("L3").Value = 334
$ws.Range("L6").Value = 268
$ws.Range("L7").Value = 950

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("L3").Value = 99
$ws.Range("L7").Value = 296

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K5").Value = 35
$ws.Range("K7").Value = 909

$ws = $wb.Worksheets.Item('New City')
$ws.Range("L2").Value = 152
$ws.Range("L7").Value = 415

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("L3").Value = 146
$ws.Range("L7").Value = 364

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("L2").Value = 84
$ws.Range("L7").Value = 210

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("L2").Value = 208
$ws.Range("L7").Value = 729

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("L3").Value = 114
$ws.Range("L6").Value = 217
$ws.Range("L7").Value = 456

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("L2").Value = 358
$ws.Range("L4").Value = 64
$ws.Range("L7").Value = 1180

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("L2").Value = 41
$ws.Range("L7").Value = 275

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("L2").Value = 207
$ws.Range("L6").Value = 156
$ws.Range("L7").Value = 574

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("L3").Value = 228
$ws.Range("L6").Value = 190
$ws.Range("L7").Value = 666

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("L4").Value = 21
$ws.Range("L6").Value = 63

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range("L6").Value = 14
$ws.Range("L7").Value = 64

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("L2").Value = 181
$ws.Range("L7").Value = 584

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("L3").Value = 185
$ws.Range("L5").Value = 6
$ws.Range("L7").Value = 532

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("L3").Value = 61
$ws.Range("L7").Value = 257

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("L2").Value = 42
$ws.Range("L7").Value = 127

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("L6").Value = 36
$ws.Range("L7").Value = 176

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("L6").Value = 57
$ws.Range("L7").Value = 182

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("L2").Value = 74
$ws.Range("L4").Value = 20
$ws.Range("L7").Value = 223

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("L2").Value = 82
$ws.Range("L6").Value = 55

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("L6").Value = 37
$ws.Range("L7").Value = 143

$ws = $wb.Worksheets.Item('Sheffield & DePaul')
$ws.Range("L6").Value = 9
$ws.Range("L7").Value = 30

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("L6").Value = 32
$ws.Range("L7").Value = 141

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range("L3").Value = 15
$ws.Range("L7").Value = 49
